{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.trim() === \"Docente(s) Respons\u00e1vel(eis)\") {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find target paragraph 'Docente(s) Respons\u00e1vel(eis)'\");\n}\n\nconst newPara = target.insertParagraph(\"4893449 - D\u00e9bora Souza Alvim\", \"After\");\nnewPara.style = \"List Bullet\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"Docente(s) Respons\u00e1vel(eis)\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find target paragraph 'Docente(s) Respons\u00e1vel(eis)'\"\n}\n\n$newRange = $target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.Text = \"4893449 - D\u00e9bora Souza Alvim\"\n$newPara.Style = $d.Styles.Item(\"List Bullet\")\n"}
